$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pairwise_comp")

$newText = "Enter pairwise comparisons in the white cells of the table or numerical data in the green cells. For the Direct Values column, if the smallest value is best, invert the value before entering it (e.g., `$10 as =1/10) ."

$ws.Range("A2").Value = $newText
$ws.Range("A10").Value = $newText
$ws.Range("A18").Value = $newText
$ws.Range("A26").Value = $newText
$ws.Range("A34").Value = $newText
$ws.Range("A42").Value = $newText
